$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G1").Value = "hello"
$v = $ws.Range("G1").Value2
Write-Host "G1:" $v
$ws.Cells.Item(1,8).Value = 42
$v2 = $ws.Cells.Item(1,8).Value2
Write-Host "H1:" $v2
